$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.566.31'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').Value = '3.477.61'
$ws.Range('E3').Value = '  -1.12%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'591.14"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.68%  '
$ws.Range('D6').Value = "'179.34"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.89%  '
$ws.Range('D7').Value = "'0.614"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.20%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = '3.476.44'
$ws.Range('E9').Value = '  -1.14%  '
$ws.Range('E10').Value = '  -2.24%  '
$ws.Range('D11').Value = "'6.98"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.66%  '
$ws.Range('E12').Value = '  -2.89%  '
$ws.Range('D13').Value = '4.084.65'
$ws.Range('E13').Value = '  -1.05%  '
$ws.Range('D14').Value = "'32.13"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.31%  '
$ws.Range('E15').Value = '  -2.56%  '
$ws.Range('D16').Value = '67.552.44'
$ws.Range('E16').Value = '  -0.48%  '
$ws.Range('E17').Value = '  -2.57%  '
$ws.Range('D18').Value = '3.474.13'
$ws.Range('E18').Value = '  -1.37%  '
$ws.Range('E19').Value = '  -3.71%  '
$ws.Range('D20').Value = "'14.05"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.92%  '
$ws.Range('D21').Value = "'385.36"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.77%  '
$ws.Range('D22').Value = "'7.90"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.18%  '
$ws.Range('D23').Value = "'5.81"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.38%  '
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').Value = "'72.09"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.30%  '
$ws.Range('D26').Value = "'0.535"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.69%  '
$ws.Range('E27').Value = '  -1.02%  '
$ws.Range('E28').Value = '  -3.86%  '
$ws.Range('E29').Value = '  -1.78%  '
$ws.Range('D30').Value = "'1.00"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.30%  '
$ws.Range('E31').Value = '  -4.34%  '
$ws.Range('D32').Value = "'24.45"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.17%  '
$ws.Range('E33').Value = '  -2.21%  '
$ws.Range('E34').Value = '  -4.89%  '
$ws.Range('E35').Value = '  -3.03%  '
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('E37').Value = '  -4.87%  '
$ws.Range('D38').Value = "'160.22"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.89%  '
$ws.Range('E39').Value = '  +0.50%  '
$ws.Range('D40').Value = "'27.76"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.06%  '
$ws.Range('D41').Value = "'1.86"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.60%  '
$ws.Range('E42').Value = '  -4.57%  '
$ws.Range('E43').Value = '  -5.11%  '
$ws.Range('E44').Value = '  -3.90%  '
$ws.Range('D45').Value = "'0.0707"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.92%  '
$ws.Range('D46').Value = '2.722.70'
$ws.Range('E46').Value = '  -6.36%  '
$ws.Range('D47').Value = "'25.82"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.21%  '
$ws.Range('D48').Value = "'41.56"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.94%  '
$ws.Range('E49').Value = '  -2.77%  '
$ws.Range('D50').Value = "'327.26"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.87%  '
$ws.Range('E51').Value = '  -3.10%  '
